# Update the "取得日時" (retrieved at) timestamp in column A for all data
# rows on the active (first) worksheet from 2025-10-09 01:16:19 to
# 2025-10-09 01:43:07, leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2025-10-09 01:16:19"
$newTimestamp = "2025-10-09 01:43:07"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
